$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 53
$ws.Range("I2").Value = 152
$ws.Range("J2").Value = 725
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 209
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = 123
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 12
$ws.Range("S2").Value = 80
$ws.Range("T2").Value = 126
$ws.Range("U2").Value = 14
$ws.Range("V2").Value = 1136
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1220
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 7
